# "Small changes to Load Balancing chapter and other"
#
# The tracking list marks, per doc chapter (rows), whether it has been
# Completed / Converted / Checked (columns B/C/D) using a green-filled,
# centered "x". The "load-balancing.md" row (row 18) was previously blank
# in all three columns; mark it fully done, matching the style already
# used by other completed rows (e.g. row 5: B5:D5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

$doneStyleSource = $ws.Range("B5:D5")
$target = $ws.Range("B18:D18")

$target.Value = "x"
$target.Interior.Color = $doneStyleSource.Interior.Color
$target.HorizontalAlignment = $doneStyleSource.HorizontalAlignment

# Leave a clean, single-cell selection instead of the stale B22:D22
# multi-cell selection that had been left over from a prior edit.
$ws.Range("A1").Select()
